# Add a new "Croatia" Test Data sheet after "Spain", modeled on the
# existing "Italy" sheet (same layout/styles/column widths), then fill in
# the Croatia-specific market name and NGC ticket reference, matching the
# "Added Test Data for Croatia Market" commit.

$wb = $excel.ActiveWorkbook

$italy = $wb.Worksheets.Item("Italy")
$spain = $wb.Worksheets.Item("Spain")

# Duplicate the Italy sheet (keeps column widths/styles/merged cells) and
# drop the copy right after Spain, i.e. at the end of the tab strip.
$italy.Copy($null, $spain)

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Croatia"

# Write the NGC reference before the market name so the shared-strings
# table grows in the same order as the source edit.
$newSheet.Range("B4").Value = "NGC-3139/T2488"
$newSheet.Range("B2").Value = "Croatia Market"

# Make the new sheet the active one, with B7 selected (matches the saved
# view state of the new sheet in the workbook).
$newSheet.Activate()
[void]$newSheet.Range("B7").Select()
